$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.693.37"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "1.895.38"
$ws.Range("E3").Value = "  +1.18%  "

$ws.Range("E4").Value = "  -1.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4887"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3796"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07332"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9147"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.56"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.31%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.933.03"
$ws.Range("E12").Value = "  +3.09%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07686"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.478"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.612"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.07"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008788"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("E19").Value = "  -1.09%  "

$ws.Range("D20").Value = "27.724.10"
$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("E21").Value = "  -2.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.127"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").Value = "2.137.52"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.75"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.906"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.52"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.158"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.84"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.906"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08912"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.193"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.223"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7660"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.646"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02035"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.529"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.096"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05285"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5489"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.983"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.919"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.526"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1518"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.84"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.25%  "

$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4801"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.03%  "

$ws.Range("E48").Value = "  -1.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.634"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.53"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06048"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.08%  "
